$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.22224168450631
$ws.Range("C2").Value = 10.71512280286257
$ws.Range("D2").Value = 6.366890670964553
$ws.Range("E2").Value = 13.22350241928203
$ws.Range("F2").Value = 30.51464489457093
$ws.Range("K2").Value = 8.610941441859012
$ws.Range("L2").Value = 9.969231619326704
$ws.Range("M2").Value = 14.07003342166148
$ws.Range("O2").Value = 27.37728486233314

$ws.Range("B3").Value = 12.01219342924657
$ws.Range("C3").Value = 10.71862745352171
$ws.Range("D3").Value = 6.331521204465656
$ws.Range("E3").Value = 13.25389325277787
$ws.Range("F3").Value = 30.5562499064167
$ws.Range("K3").Value = 8.448034244293428
$ws.Range("L3").Value = 9.977100882347298
$ws.Range("M3").Value = 14.0425839426659
$ws.Range("O3").Value = 27.45032074583878

$ws.Range("B4").Value = 11.88372986712437
$ws.Range("C4").Value = 10.72113088122725
$ws.Range("D4").Value = 6.309310850992324
$ws.Range("E4").Value = 13.27431673644092
$ws.Range("F4").Value = 30.58872595969697
$ws.Range("K4").Value = 8.347819480564226
$ws.Range("L4").Value = 9.983261192141319
$ws.Range("M4").Value = 14.02779152086773
$ws.Range("O4").Value = 27.5003072277037

$ws.Range("B5").Value = 11.83157727433641
$ws.Range("C5").Value = 10.72223984542581
$ws.Range("D5").Value = 6.300138251206045
$ws.Range("E5").Value = 13.28308310912166
$ws.Range("F5").Value = 30.60370055539411
$ws.Range("K5").Value = 8.306988566370709
$ws.Range("L5").Value = 9.986106150724966
$ws.Range("M5").Value = 14.02228649506477
$ws.Range("O5").Value = 27.5219683442886

$ws.Range("B6").Value = 11.82293134906925
$ws.Range("C6").Value = 10.72242936282233
$ws.Range("D6").Value = 6.298607866714817
$ws.Range("E6").Value = 13.28456555939744
$ws.Range("F6").Value = 30.60629210701279
$ws.Range("K6").Value = 8.300210699100059
$ws.Range("L6").Value = 9.986598775497566
$ws.Range("M6").Value = 14.0214041049708
$ws.Range("O6").Value = 27.52564308733054

$ws.Range("B7").Value = 11.88302562718413
$ws.Range("C7").Value = 10.72114547704535
$ws.Range("D7").Value = 6.309187635750181
$ws.Range("E7").Value = 13.27443316625394
$ws.Range("F7").Value = 30.58892086934986
$ws.Range("K7").Value = 8.347268717140512
$ws.Range("L7").Value = 9.983298204909795
$ws.Range("M7").Value = 14.02771515474154
$ws.Range("O7").Value = 27.50059413150345

$ws.Range("B8").Value = 12.14975354108302
$ws.Range("C8").Value = 10.71625847630743
$ws.Range("D8").Value = 6.354799082283145
$ws.Range("E8").Value = 13.23361532287112
$ws.Range("F8").Value = 30.52755064520947
$ws.Range("K8").Value = 8.554842898146537
$ws.Range("L8").Value = 9.971669521759262
$ws.Range("M8").Value = 14.06014379282194
$ws.Range("O8").Value = 27.40139932169662

$ws.Range("B9").Value = 12.67380393761587
$ws.Range("C9").Value = 10.7094475448984
$ws.Range("D9").Value = 6.440227282138985
$ws.Range("E9").Value = 13.16755344902248
$ws.Range("F9").Value = 30.46227634945155
$ws.Range("K9").Value = 8.958046491994672
$ws.Range("L9").Value = 9.959383263496688
$ws.Range("M9").Value = 14.13987959463995
$ws.Range("O9").Value = 27.24775697367454

$ws.Range("B10").Value = 13.05548486233765
$ws.Range("C10").Value = 10.70611140046294
$ws.Range("D10").Value = 6.500411454563328
$ws.Range("E10").Value = 13.12752766114864
$ws.Range("F10").Value = 30.44797197569103
$ws.Range("K10").Value = 9.24890929141139
$ws.Range("L10").Value = 9.956733812941389
$ws.Range("M10").Value = 14.2079949380807
$ws.Range("O10").Value = 27.15989728792775

$ws.Range("B11").Value = 13.2276011006313
$ws.Range("C11").Value = 10.70495117608615
$ws.Range("D11").Value = 6.527203043353349
$ws.Range("E11").Value = 13.11116354257828
$ws.Range("F11").Value = 30.44877527234193
$ws.Range("K11").Value = 9.379465579040657
$ws.Range("L11").Value = 9.956904597120223
$ws.Range("M11").Value = 14.24097991776724
$ws.Range("O11").Value = 27.12537982568984

$ws.Range("B12").Value = 13.29249677107257
$ws.Range("C12").Value = 10.70456283101706
$ws.Range("D12").Value = 6.537261627118307
$ws.Range("E12").Value = 13.10523173902085
$ws.Range("F12").Value = 30.45012969705645
$ws.Range("K12").Value = 9.428604362175424
$ws.Range("L12").Value = 9.957166268275031
$ws.Range("M12").Value = 14.25375122667106
$ws.Range("O12").Value = 27.1130940432745

$ws.Range("B13").Value = 13.2785338049754
$ws.Range("C13").Value = 10.70464420569414
$ws.Range("D13").Value = 6.535099231062889
$ws.Range("E13").Value = 13.10649747955483
$ws.Range("F13").Value = 30.44979131165884
$ws.Range("K13").Value = 9.418035513045254
$ws.Range("L13").Value = 9.957101165821765
$ws.Range("M13").Value = 14.25098832859551
$ws.Range("O13").Value = 27.11570506113389

$ws.Range("B14").Value = 13.23294609163702
$ws.Range("C14").Value = 10.70491820666788
$ws.Range("D14").Value = 6.528032322943879
$ws.Range("E14").Value = 13.11067022111966
$ws.Range("F14").Value = 30.44886566173036
$ws.Range("K14").Value = 9.383514530298305
$ws.Range("L14").Value = 9.956922182075727
$ws.Range("M14").Value = 14.24202503739328
$ws.Range("O14").Value = 27.12435332187504

$ws.Range("B15").Value = 13.20498387026173
$ws.Range("C15").Value = 10.70509267127797
$ws.Range("D15").Value = 6.523692259626499
$ws.Range("E15").Value = 13.11326064149972
$ws.Range("F15").Value = 30.4484354038368
$ws.Range("K15").Value = 9.362329002533752
$ws.Range("L15").Value = 9.956838177261783
$ws.Range("M15").Value = 14.23657109862348
$ws.Range("O15").Value = 27.12975293272815

$ws.Range("B16").Value = 13.04420072950749
$ws.Range("C16").Value = 10.70619438257934
$ws.Range("D16").Value = 6.498648547486691
$ws.Range("E16").Value = 13.12863417991723
$ws.Range("F16").Value = 30.44806651424295
$ws.Range("K16").Value = 9.240337650294979
$ws.Range("L16").Value = 9.956750268161379
$ws.Range("M16").Value = 14.20587893920752
$ws.Range("O16").Value = 27.16226289018365

$ws.Range("B17").Value = 12.94513198690432
$ws.Range("C17").Value = 10.70696151038303
$ws.Range("D17").Value = 6.483133058291475
$ws.Range("E17").Value = 13.13853745261183
$ws.Range("F17").Value = 30.44971210955036
$ws.Range("K17").Value = 9.165015078039831
$ws.Range("L17").Value = 9.9570481761703
$ws.Range("M17").Value = 14.18755751333149
$ws.Range("O17").Value = 27.18360376307099

$ws.Range("B18").Value = 12.88801153824383
$ws.Range("C18").Value = 10.70743640568601
$ws.Range("D18").Value = 6.474153901935159
$ws.Range("E18").Value = 13.14440711499496
$ws.Range("F18").Value = 30.45134680995911
$ws.Range("K18").Value = 9.121528648583514
$ws.Range("L18").Value = 9.957349100582135
$ws.Range("M18").Value = 14.17720809553748
$ws.Range("O18").Value = 27.19639140283581

$ws.Range("B19").Value = 12.86864969392421
$ws.Range("C19").Value = 10.70760299067989
$ws.Range("D19").Value = 6.471104333527776
$ws.Range("E19").Value = 13.14642429717945
$ws.Range("F19").Value = 30.45201850852557
$ws.Range("K19").Value = 9.106778396746124
$ws.Range("L19").Value = 9.957473268534606
$ws.Range("M19").Value = 14.17373655275498
$ws.Range("O19").Value = 27.2008091272726

$ws.Range("B20").Value = 12.95569285131799
$ws.Range("C20").Value = 10.70687636706675
$ws.Range("D20").Value = 6.484790424824159
$ws.Range("E20").Value = 13.13746527075372
$ws.Range("F20").Value = 30.44946571215575
$ws.Range("K20").Value = 9.173050501784624
$ws.Range("L20").Value = 9.957003058855051
$ws.Range("M20").Value = 14.18948839141435
$ws.Range("O20").Value = 27.18127889273415

$ws.Range("B21").Value = 13.24634442613794
$ws.Range("C21").Value = 10.70483634471428
$ws.Range("D21").Value = 6.530110420705995
$ws.Range("E21").Value = 13.10943739729733
$ws.Range("F21").Value = 30.44910905606687
$ws.Range("K21").Value = 9.393662685481443
$ws.Range("L21").Value = 9.956969414772136
$ws.Range("M21").Value = 14.24465021355475
$ws.Range("O21").Value = 27.12179179611961

$ws.Range("B22").Value = 13.43463677068195
$ws.Range("C22").Value = 10.70380019956095
$ws.Range("D22").Value = 6.559222540147363
$ws.Range("E22").Value = 13.09266368524732
$ws.Range("F22").Value = 30.45499654809165
$ws.Range("K22").Value = 9.53607662314365
$ws.Range("L22").Value = 9.958095229306425
$ws.Range("M22").Value = 14.28233389400162
$ws.Range("O22").Value = 27.08749104046269

$ws.Range("B23").Value = 13.33431417610005
$ws.Range("C23").Value = 10.7043261502395
$ws.Range("D23").Value = 6.543732064599293
$ws.Range("E23").Value = 13.10147491635822
$ws.Range("F23").Value = 30.45129476660921
$ws.Range("K23").Value = 9.460244357410911
$ws.Range("L23").Value = 9.957389644434514
$ws.Range("M23").Value = 14.26207439208164
$ws.Range("O23").Value = 27.105378711139

$ws.Range("B24").Value = 12.95091879205743
$ws.Range("C24").Value = 10.70691475485001
$ws.Range("D24").Value = 6.484041313288484
$ws.Range("E24").Value = 13.13794945550361
$ws.Range("F24").Value = 30.44957496331157
$ws.Range("K24").Value = 9.169418251477595
$ws.Range("L24").Value = 9.957023052474675
$ws.Range("M24").Value = 14.18861486835734
$ws.Range("O24").Value = 27.18232835164717

$ws.Range("B25").Value = 12.53234256692108
$ws.Range("C25").Value = 10.71099558358701
$ws.Range("D25").Value = 6.41756027395244
$ws.Range("E25").Value = 13.18392945352999
$ws.Range("F25").Value = 30.47402808996067
$ws.Range("K25").Value = 8.849701370354172
$ws.Range("L25").Value = 9.961584333315461
$ws.Range("M25").Value = 14.11661110055485
$ws.Range("O25").Value = 27.2849340492647
